$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 984.1111
$ws.Range("I9").Value = 244.25
$ws.Range("J9").Value = 1576
$ws.Range("K9").Value = 244.25
$ws.Range("L9").Value = 1576
$ws.Range("M9").Value = -75.25
$ws.Range("N9").Value = -1914
$ws.Range("H41").Value = 1299.7142
$ws.Range("I41").Value = 1149.6666
$ws.Range("J41").Value = 2200
$ws.Range("K41").Value = 1149.6666
$ws.Range("L41").Value = 2200
$ws.Range("M41").Value = -709.6666
$ws.Range("N41").Value = -3080
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H81").Value = 296326.66
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 296326.66
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H98").Value = 3833.3333
$ws.Range("I98").Value = 4375
$ws.Range("J98").Value = 2750
$ws.Range("K98").Value = 4375
$ws.Range("L98").Value = 2750
$ws.Range("M98").Value = -2877
$ws.Range("N98").Value = -5746
$ws.Range("H122").Value = 3833.3333
$ws.Range("I122").Value = 4375
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 13125
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -10675
$ws.Range("N122").Value = -13150
$ws.Range("H127").Value = 357948.56
$ws.Range("I127").Value = 455053.28
$ws.Range("K127").Value = 1365159.84
$ws.Range("M127").Value = -1360199.84
$ws.Range("H138").Value = 3046.6667
$ws.Range("I138").Value = 3707.4707
$ws.Range("J138").Value = 2455.4211
$ws.Range("K138").Value = 11122.4121
$ws.Range("L138").Value = 7366.263300000001
$ws.Range("M138").Value = -5982.4121
$ws.Range("N138").Value = -17646.2633

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H5").Value = 3300
$ws.Range("I5").Value = 3920
$ws.Range("K5").Value = 3920
$ws.Range("M5").Value = -3808
$ws.Range("H36").Value = 3440.625
$ws.Range("I36").Value = 3440.625
$ws.Range("K36").Value = 3440.625
$ws.Range("M36").Value = -3094.625
$ws.Range("H61").Value = 19967
$ws.Range("I61").Value = 5863.875
$ws.Range("K61").Value = 5863.875
$ws.Range("M61").Value = -5651.875
$ws.Range("H97").Value = 3548.75
$ws.Range("I97").Value = 1461.7273
$ws.Range("K97").Value = 1461.7273
$ws.Range("M97").Value = -965.7273
$ws.Range("H136").Value = 19967
$ws.Range("I136").Value = 5863.875
$ws.Range("K136").Value = 17591.625
$ws.Range("M136").Value = -15041.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 71569.60000000001
$ws.Range("J2").Value = 71569.60000000001
$ws.Range("L2").Value = 71569.60000000001
$ws.Range("N2").Value = -71795.60000000001
$ws.Range("H4").Value = 3300
$ws.Range("I4").Value = 3920
$ws.Range("K4").Value = 3920
$ws.Range("M4").Value = -3805
$ws.Range("H33").Value = 10584.4
$ws.Range("I33").Value = 2949
$ws.Range("J33").Value = 15674.667
$ws.Range("K33").Value = 2949
$ws.Range("L33").Value = 15674.667
$ws.Range("M33").Value = -2613
$ws.Range("N33").Value = -16346.667
$ws.Range("H80").Value = 1816.4166
$ws.Range("J80").Value = 2049.1
$ws.Range("L80").Value = 2049.1
$ws.Range("N80").Value = -4045.1
$ws.Range("H83").Value = 1816.4166
$ws.Range("J83").Value = 2049.1
$ws.Range("L83").Value = 10245.5
$ws.Range("N83").Value = -20229.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 16529.482
$ws.Range("I58").Value = 13750.4
$ws.Range("J58").Value = 17108.459
$ws.Range("K58").Value = 13750.4
$ws.Range("L58").Value = 17108.459
$ws.Range("M58").Value = -13547.4
$ws.Range("N58").Value = -17514.459
$ws.Range("H94").Value = 4902.3
$ws.Range("J94").Value = 5689
$ws.Range("L94").Value = 5689
$ws.Range("N94").Value = -6591
$ws.Range("H136").Value = 16529.482
$ws.Range("I136").Value = 13750.4
$ws.Range("J136").Value = 17108.459
$ws.Range("K136").Value = 41251.2
$ws.Range("L136").Value = 51325.37699999999
$ws.Range("M136").Value = -38701.2
$ws.Range("N136").Value = -56425.37699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 43529244
$ws.Range("I4").Value = 47395536
$ws.Range("K4").Value = 142186608
$ws.Range("M4").Value = -142186496
$ws.Range("H92").Value = 873.6667
$ws.Range("I92").Value = 1275
$ws.Range("J92").Value = 673
$ws.Range("K92").Value = 3825
$ws.Range("L92").Value = 2019
$ws.Range("M92").Value = -2577
$ws.Range("N92").Value = -4515
$ws.Range("H97").Value = 8931.25
$ws.Range("I97").Value = 487.5
$ws.Range("J97").Value = 17375
$ws.Range("K97").Value = 1462.5
$ws.Range("L97").Value = 52125
$ws.Range("M97").Value = -966.5
$ws.Range("N97").Value = -53117
$ws.Range("H98").Value = 5261.25
$ws.Range("J98").Value = 5955.7144
$ws.Range("L98").Value = 17867.1432
$ws.Range("N98").Value = -20863.1432
$ws.Range("H117").Value = 1842.7142
$ws.Range("J117").Value = 2099.8333
$ws.Range("L117").Value = 6299.499899999999
$ws.Range("N117").Value = -13183.4999
$ws.Range("H121").Value = 2068.5715
$ws.Range("I121").Value = 200
$ws.Range("J121").Value = 2380
$ws.Range("K121").Value = 600
$ws.Range("L121").Value = 7140
$ws.Range("M121").Value = 710
$ws.Range("N121").Value = -9760
$ws.Range("H122").Value = 11329894
$ws.Range("I122").Value = 18687766
$ws.Range("K122").Value = 168189894
$ws.Range("M122").Value = -168187444
$ws.Range("H129").Value = 2552
$ws.Range("I129").Value = 2216
$ws.Range("J129").Value = 3000
$ws.Range("K129").Value = 6648
$ws.Range("L129").Value = 9000
$ws.Range("M129").Value = -1648
$ws.Range("N129").Value = -19000
$ws.Range("H131").Value = 1447.13
$ws.Range("I131").Value = 942.6667
$ws.Range("J131").Value = 1479.3298
$ws.Range("K131").Value = 2828.0001
$ws.Range("L131").Value = 4437.9894
$ws.Range("M131").Value = 2211.9999
$ws.Range("N131").Value = -14517.9894

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1429114
$ws.Range("I3").Value = 2000339.6
$ws.Range("J3").Value = 1050
$ws.Range("K3").Value = 2000339.6
$ws.Range("L3").Value = 1050
$ws.Range("M3").Value = -2000223.6
$ws.Range("N3").Value = -1282
$ws.Range("H14").Value = 1250152.1
$ws.Range("I14").Value = 2000140.4
$ws.Range("K14").Value = 2000140.4
$ws.Range("M14").Value = -1999972.4
$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H52").Value = 28889.525
$ws.Range("J52").Value = 28889.525
$ws.Range("L52").Value = 28889.525
$ws.Range("N52").Value = -29407.525

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3467.5186
$ws.Range("I61").Value = 2137.0625
$ws.Range("J61").Value = 5402.727
$ws.Range("K61").Value = 2137.0625
$ws.Range("L61").Value = 5402.727
$ws.Range("M61").Value = -1935.0625
$ws.Range("N61").Value = -5806.727
$ws.Range("H113").Value = 3467.5186
$ws.Range("I113").Value = 2137.0625
$ws.Range("J113").Value = 5402.727
$ws.Range("K113").Value = 2137.0625
$ws.Range("L113").Value = 5402.727
$ws.Range("M113").Value = 32.9375
$ws.Range("N113").Value = -9742.726999999999
$ws.Range("H132").Value = 2362860
$ws.Range("I132").Value = 4031.7
$ws.Range("J132").Value = 5732615
$ws.Range("K132").Value = 12095.1
$ws.Range("L132").Value = 17197845
$ws.Range("M132").Value = -9565.099999999999
$ws.Range("N132").Value = -17202905
$ws.Range("H133").Value = 61990
$ws.Range("J133").Value = 61990
$ws.Range("L133").Value = 61990
$ws.Range("N133").Value = -67050

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 13500.5
$ws.Range("I3").Value = 1250
$ws.Range("K3").Value = 1250
$ws.Range("M3").Value = -1136
$ws.Range("H100").Value = 950.1
$ws.Range("I100").Value = 936.5
$ws.Range("K100").Value = 1873
$ws.Range("M100").Value = -1332
$ws.Range("H101").Value = 28464.572
$ws.Range("J101").Value = 26542
$ws.Range("L101").Value = 26542
$ws.Range("N101").Value = -33032
$ws.Range("H107").Value = 2357.9119
$ws.Range("I107").Value = 862.087
$ws.Range("J107").Value = 5485.5454
$ws.Range("K107").Value = 2586.261
$ws.Range("L107").Value = 16456.6362
$ws.Range("M107").Value = -666.261
$ws.Range("N107").Value = -20296.6362
